# Revert "Release for 2nd version"
#
# The release commit appended two extra slides ("Dinh nghia" / "Phan loai",
# slide4.xml + slide5.xml) to a 3-slide deck. Reverting it removes those two
# trailing slides again, leaving slides 1-3 ("KY NGHE PHAN MEM NANG CAO",
# "YEU CAU PHAN MEM" / "CHUONG 1", "Thuat ngu thuong dung") untouched.

$p = $ppt.ActivePresentation

# Remove the two trailing slides that the release commit had added.
# Walk from the end so indices of the slides we keep never shift under us.
for ($i = $p.Slides.Count; $i -ge 4; $i--) {
    $p.Slides.Item($i).Delete()
}
